$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 54
$ws1.Range("F6").Value = 65
$ws1.Range("F8").Value = 3871
$ws1.Range("F10").Value = 4560
$ws1.Range("F12").Value = 1149
$ws1.Range("F13").Value = 70

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 54
$ws4.Range("F6").Value = 65
$ws4.Range("F9").Value = 3871
$ws4.Range("F11").Value = 4560
$ws4.Range("F13").Value = 1149
$ws4.Range("F14").Value = 70
